$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $cell = $ws.Range($rangeAddress)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.477.74"

# Row 3 - Ethereum
Set-TextValue "D3" "2.469.61"
$ws.Range("E3").Value = "  +2.11%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
Set-TextValue "D5" "575.79"
$ws.Range("E5").Value = "  +2.33%  "

# Row 6 - Solana
Set-TextValue "D6" "148.07"
$ws.Range("E6").Value = "  +3.57%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.09%  "

# Row 8 - XRP
Set-TextValue "D8" "0.540"
$ws.Range("E8").Value = "  +1.70%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +4.24%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  +0.64%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.33"
$ws.Range("E11").Value = "  +2.76%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.362"
$ws.Range("E12").Value = "  +3.84%  "

# Row 13 - Avalanche
Set-TextValue "D13" "27.23"
$ws.Range("E13").Value = "  +3.95%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +6.37%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.911.17"
$ws.Range("E15").Value = "  +2.22%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "63.472.18"
$ws.Range("E16").Value = "  +2.45%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.479.27"
$ws.Range("E17").Value = "  +2.79%  "

# Row 18 - Chainlink
Set-TextValue "D18" "11.54"
$ws.Range("E18").Value = "  +1.78%  "

# Row 19 - Uniswap
Set-TextValue "D19" "7.29"
$ws.Range("E19").Value = "  +6.91%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +2.48%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "328.56"
$ws.Range("E21").Value = "  +1.55%  "

# Row 22 - Dai
Set-TextValue "D22" "0.998"
$ws.Range("E22").Value = "  -0.17%  "

# Row 23 - SuiNetwork
$ws.Range("E23").Value = "  +10.88%  "

# Row 24 - Litecoin
Set-TextValue "D24" "67.40"
$ws.Range("E24").Value = "  +1.01%  "

# Row 25 - Bittensor
Set-TextValue "D25" "627.91"
$ws.Range("E25").Value = "  +13.67%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +12.97%  "

# Row 27 - Aptos
Set-TextValue "D27" "8.74"
$ws.Range("E27").Value = "  -0.32%  "

# Row 28 - WrappedeETH
$ws.Range("E28").Value = "  +2.35%  "

# Row 29 - Fetch.AI
$ws.Range("E29").Value = "  +9.45%  "

# Row 30 & 31 - swapped: InternetComputer(DFINITY) and Binance-PegBSC-USD
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "8.45"
$ws.Range("E31").Value = "  +2.76%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  -1.35%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +2.54%  "

# Row 34 - NEARProtocol
Set-TextValue "D34" "5.18"
$ws.Range("E34").Value = "  +9.58%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +3.43%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  -0.14%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("E37").Value = "  +1.96%  "

# Row 38 - RenderToken
Set-TextValue "D38" "5.52"
$ws.Range("E38").Value = "  +1.64%  "

# Row 39 - EthereumClassic
Set-TextValue "D39" "18.95"
$ws.Range("E39").Value = "  +2.13%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +2.02%  "

# Row 41 - Monero
Set-TextValue "D41" "147.17"
$ws.Range("E41").Value = "  -3.50%  "

# Row 42 - dogwifhat
Set-TextValue "D42" "2.67"
$ws.Range("E42").Value = "  +20.13%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  -0.02%  "

# Row 44 - Aave
Set-TextValue "D44" "150.25"
$ws.Range("E44").Value = "  +2.06%  "

# Row 45 - Filecoin
Set-TextValue "D45" "3.76"
$ws.Range("E45").Value = "  +3.44%  "

# Row 46 - Hedera
Set-TextValue "D46" "0.0549"
$ws.Range("E46").Value = "  +4.17%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "21.15"
$ws.Range("E47").Value = "  +6.65%  "

# Row 48 - Mantle
Set-TextValue "D48" "0.609"
$ws.Range("E48").Value = "  +2.50%  "

# Row 49 - VeChain
Set-TextValue "D49" "0.0240"
$ws.Range("E49").Value = "  +5.65%  "

# Row 50 - Stellar
Set-TextValue "D50" "0.0927"
$ws.Range("E50").Value = "  +0.80%  "

# Row 51 - ONDO
Set-TextValue "D51" "0.746"
$ws.Range("E51").Value = "  +4.77%  "
